# Generate Report for Handoff
#
# A new localization request (174186a9-1bfd-4d29-ab0d-7a037b9a15da) supersedes the
# previous one (f0d15859-f377-4e16-94c0-329d2d4de750): refresh every file-name /
# hash reference, bump the "Latest HO Xliff Generate Date" and the per-language
# "Latest Handoff Datetime", and since the freshly generated xliff hasn't been
# handed back yet, clear "Latest Target File" / "Latest Handback File" and reset
# "Latest Handback DateTime" to the zero date.

$wb = $excel.ActiveWorkbook

$oldId = "f0d15859-f377-4e16-94c0-329d2d4de750"
$newId = "174186a9-1bfd-4d29-ab0d-7a037b9a15da"
$oldZhHash = "5a7cd9757f250c791862152d2b7ed58834670876"
$newZhHash = "dfa032418414a1599b2ecdfd20e1130f222e5be7"
$oldDeHash = "5a7cd9757f250c791862152d2b7ed58834670876"
$newDeHash = "dfa032418414a1599b2ecdfd20e1130f222e5be7"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-22 13:03:02"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newId.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = "$newId.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-22 13:02:55"

# drop the I2 hyperlink (the file has not been handed back yet)
foreach ($hl in @($wsZh.Hyperlinks)) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}
$wsZh.Range("I2").Value = "'"
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = "'"
$wsZh.Range("J2").Style = "Normal"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newId.md"
    }
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = "$newId.$newDeHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-22 13:03:02"

foreach ($hl in @($wsDe.Hyperlinks)) {
    if ($hl.Range.Address() -eq '$I$2') {
        $hl.Delete()
    }
}
$wsDe.Range("I2").Value = "'"
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = "'"
$wsDe.Range("J2").Style = "Normal"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = "$newId.md"
    }
}

Write-Host "Report regenerated for handoff $newId"
